$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 4875.4546
$ws.Range("I18").Value = 2338.2
$ws.Range("J18").Value = 10312.429
$ws.Range("K18").Value = 2338.2
$ws.Range("L18").Value = 10312.429
$ws.Range("M18").Value = -2054.2
$ws.Range("N18").Value = -10880.429
$ws.Range("H19").Value = 1117.2142
$ws.Range("J19").Value = 1460.3334
$ws.Range("L19").Value = 1460.3334
$ws.Range("N19").Value = -1810.3334
$ws.Range("H62").Value = 2000
$ws.Range("I62").Value = 2000
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 2000
$ws.Range("L62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -1376
$ws.Range("H65").Value = 2000
$ws.Range("I65").Value = 2000
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 10000
$ws.Range("L65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = -6880
$ws.Range("H80").Value = 4101
$ws.Range("I80").Value = 4101
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 12303
$ws.Range("L80").Value = 0
$ws.Range("M80").ClearContents()
$ws.Range("N80").Value = -11305
$ws.Range("H83").Value = 4101
$ws.Range("I83").Value = 4101
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 36909
$ws.Range("L83").Value = 0
$ws.Range("M83").ClearContents()
$ws.Range("N83").Value = -31917
$ws.Range("H97").Value = 1092.6666
$ws.Range("J97").Value = 1195
$ws.Range("L97").Value = 3585
$ws.Range("N97").Value = -4577
$ws.Range("H98").Value = 2997.6667
$ws.Range("I98").Value = 3097.5
$ws.Range("J98").Value = 2897.8333
$ws.Range("K98").Value = 3097.5
$ws.Range("L98").Value = 2897.8333
$ws.Range("M98").Value = -1599.5
$ws.Range("N98").Value = -5893.8333
$ws.Range("H122").Value = 2997.6667
$ws.Range("I122").Value = 3097.5
$ws.Range("J122").Value = 2897.8333
$ws.Range("K122").Value = 9292.5
$ws.Range("L122").Value = 8693.499899999999
$ws.Range("M122").Value = -6842.5
$ws.Range("N122").Value = -13593.4999
$ws.Range("H137").Value = 2210.3333
$ws.Range("I137").Value = 1578.8334
$ws.Range("J137").Value = 2631.3333
$ws.Range("K137").Value = 4736.5002
$ws.Range("L137").Value = 7893.999899999999
$ws.Range("M137").Value = -2186.5002
$ws.Range("N137").Value = -12993.9999
$ws.Range("H138").Value = 4529.579
$ws.Range("I138").Value = 5594.222
$ws.Range("K138").Value = 16782.666
$ws.Range("M138").Value = -11642.666
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5285.773
$ws.Range("I32").Value = 4291.579
$ws.Range("J32").Value = 11582.333
$ws.Range("K32").Value = 4291.579
$ws.Range("L32").Value = 11582.333
$ws.Range("M32").Value = -4004.579
$ws.Range("N32").Value = -12156.333
$ws.Range("H45").Value = 1800.2727
$ws.Range("I45").Value = 912
$ws.Range("J45").Value = 1889.1
$ws.Range("K45").Value = 912
$ws.Range("L45").Value = 1889.1
$ws.Range("M45").Value = -535
$ws.Range("N45").Value = -2643.1
$ws.Range("H61").Value = 3665.111
$ws.Range("I61").Value = 1741
$ws.Range("K61").Value = 1741
$ws.Range("M61").Value = -1529
$ws.Range("H110").Value = 232.33333
$ws.Range("I110").Value = 213.5
$ws.Range("K110").Value = 213.5
$ws.Range("M110").Value = 1831.5
$ws.Range("H122").Value = 2250.65
$ws.Range("I122").Value = 1808.3334
$ws.Range("J122").Value = 3577.6
$ws.Range("K122").Value = 5425.0002
$ws.Range("L122").Value = 10732.8
$ws.Range("M122").Value = -2975.0002
$ws.Range("N122").Value = -15632.8
$ws.Range("H136").Value = 3665.111
$ws.Range("I136").Value = 1741
$ws.Range("K136").Value = 5223
$ws.Range("M136").Value = -2673
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 13303.625
$ws.Range("I80").Value = 700
$ws.Range("J80").Value = 17504.834
$ws.Range("K80").Value = 700
$ws.Range("L80").Value = 17504.834
$ws.Range("M80").Value = 298
$ws.Range("N80").Value = -19500.834
$ws.Range("H83").Value = 13303.625
$ws.Range("I83").Value = 700
$ws.Range("J83").Value = 17504.834
$ws.Range("K83").Value = 3500
$ws.Range("L83").Value = 87524.17
$ws.Range("M83").Value = 1492
$ws.Range("N83").Value = -97508.17
$ws.Range("H134").Value = 5107.241
$ws.Range("I134").Value = 5266.7856
$ws.Range("K134").Value = 15800.3568
$ws.Range("M134").Value = -13265.3568
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 798.8
$ws.Range("I10").Value = 798.8
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 798.8
$ws.Range("L10").Value = 0
$ws.Range("M10").ClearContents()
$ws.Range("N10").Value = -659.8
$ws.Range("H58").Value = 3716.75
$ws.Range("I58").Value = 2831.7778
$ws.Range("J58").Value = 4854.5713
$ws.Range("K58").Value = 2831.7778
$ws.Range("L58").Value = 4854.5713
$ws.Range("M58").Value = -2628.7778
$ws.Range("N58").Value = -5260.5713
$ws.Range("H132").Value = 2804.7856
$ws.Range("I132").Value = 2066.3809
$ws.Range("K132").Value = 6199.1427
$ws.Range("M132").Value = -3669.1427
$ws.Range("H136").Value = 3716.75
$ws.Range("I136").Value = 2831.7778
$ws.Range("J136").Value = 4854.5713
$ws.Range("K136").Value = 8495.3334
$ws.Range("L136").Value = 14563.7139
$ws.Range("M136").Value = -5945.3334
$ws.Range("N136").Value = -19663.7139
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("L37").ClearContents()
$ws.Range("N37").Value = 0
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").ClearContents()
$ws.Range("N68").Value = 0
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").ClearContents()
$ws.Range("N71").Value = 0
$ws.Range("H97").Value = 350
$ws.Range("I97").Value = 200
$ws.Range("K97").Value = 600
$ws.Range("M97").Value = -104
$ws.Range("H122").Value = 1996.091
$ws.Range("J122").Value = 2095.8
$ws.Range("L122").Value = 18862.2
$ws.Range("N122").Value = -23762.2
$ws.Range("H140").Value = 2164.9106
$ws.Range("I140").Value = 1164.4762
$ws.Range("K140").Value = 3493.4286
$ws.Range("M140").Value = 1686.5714
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 3208.0833
$ws.Range("I70").Value = 2583
$ws.Range("K70").Value = 2583
$ws.Range("M70").Value = -2313
$ws.Range("H73").Value = 3208.0833
$ws.Range("I73").Value = 2583
$ws.Range("K73").Value = 2583
$ws.Range("M73").Value = -1647
$ws.Range("H102").Value = 3402.9355
$ws.Range("I102").Value = 3780.9
$ws.Range("K102").Value = 3780.9
$ws.Range("M102").Value = -2158.9
$ws.Range("H122").Value = 1691.1666
$ws.Range("I122").Value = 1691.1666
$ws.Range("K122").Value = 5073.4998
$ws.Range("M122").Value = -2623.4998
$ws.Range("H139").Value = 65000
$ws.Range("J139").Value = 65000
$ws.Range("L139").Value = 65000
$ws.Range("N139").Value = -75280
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4620.909
$ws.Range("I7").Value = 3659.1667
$ws.Range("K7").Value = 3659.1667
$ws.Range("M7").Value = -3547.1667
$ws.Range("H126").Value = 4620.909
$ws.Range("I126").Value = 3659.1667
$ws.Range("K126").Value = 10977.5001
$ws.Range("M126").Value = -8507.500100000001
$ws.Range("H136").Value = 3460.0881
$ws.Range("I136").Value = 2678.56
$ws.Range("K136").Value = 8035.68
$ws.Range("M136").Value = -5485.68
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H109").Value = 78879.5
$ws.Range("J109").Value = 78879.5
$ws.Range("L109").Value = 78879.5
$ws.Range("N109").Value = -81653.5

$wb.Save()